$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10; this shifts the existing weekly records
# (rows 10-15) down to rows 11-16, preserving all their data.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44893
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 900
$ws.Cells.Item(10, 11).Value = 13000
$ws.Cells.Item(10, 12).Value = 14000
$ws.Cells.Item(10, 13).Value = 13444
$ws.Cells.Item(10, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 1034
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"
